# Circuit_Components.xlsx update
# - Add chosen diode (SP4203-01FTG-C) and op amp (OPA172IDR) entries to the
#   components spreadsheet (ESD Protection diode + 4 extra op-amp rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$currencyFmt = '"$"#,##0.00;[Red]\-"$"#,##0.00'

# ---------------------------------------------------------------------------
# Row 11 (Driven right leg circuit op amp): price was missing, now populated.
# ---------------------------------------------------------------------------
$ws.Range("F11").Value = 2.28
$ws.Range("F11").NumberFormat = $currencyFmt

# ---------------------------------------------------------------------------
# Helper: fill in an op-amp row (Part/Value/Price) and wire up the Digikey
# hyperlink for the "OPA172IDR" part, same as the existing D11 hyperlink.
# ---------------------------------------------------------------------------
function Add-OpAmpRow($row) {
    $cell = "D" + $row
    $ws.Hyperlinks.Add($ws.Range($cell), "https://www.digikey.com.au/en/products/detail/texas-instruments/OPA172IDR/4695363", "", "", "https://www.digikey.com.au/en/products/detail/texas-instruments/OPA172IDR/4695363")
    $ws.Range($cell).Value = "OPA172IDR"
    $ws.Range($cell).Style = "Hyperlink"

    $fcell = "F" + $row
    $ws.Range($fcell).Value = 2.28
    $ws.Range($fcell).NumberFormat = $currencyFmt
}

Add-OpAmpRow 16
Add-OpAmpRow 20
Add-OpAmpRow 25
Add-OpAmpRow 28

# ---------------------------------------------------------------------------
# Row 33/34 (new "ESD Protection" section): chosen diode SP4203-01FTG-C.
# ---------------------------------------------------------------------------
$ws.Range("G33").Value = $null
$ws.Range("G33").NumberFormat = $currencyFmt

$ws.Hyperlinks.Add($ws.Range("D34"), "https://www.digikey.com.au/en/products/detail/littelfuse-inc/SP4203-01FTG-C/9828985", "", "", "https://www.digikey.com.au/en/products/detail/littelfuse-inc/SP4203-01FTG-C/9828985")
$ws.Range("D34").Value = "SP4203-01FTG-C"
$ws.Range("D34").Style = "Hyperlink"

$ws.Range("F34").Value = 0.99
$ws.Range("F34").NumberFormat = $currencyFmt

$ws.Range("G34").Formula = "=E34*F34"
$ws.Range("G34").NumberFormat = $currencyFmt

# ---------------------------------------------------------------------------
# Selection cursor, matching the saved workbook view.
# ---------------------------------------------------------------------------
$ws.Range("I12").Select()
